# Apply weekly update: insert a new data row at row 50 (shifting existing
# rows 50-156 down to 51-157) and populate the new row with the latest
# market data (date 2021-09-30 / serial 44469).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 50, shifting rows 50:156 down to 51:157
$ws.Rows.Item(50).Insert()

# Fill in the new row 50 with the latest data entry
$ws.Range("A50").Value = 9
$ws.Range("B50").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C50").Value = "Metropolitana"
$ws.Range("D50").Value = 44469
$ws.Range("E50").Value = 13
$ws.Range("F50").Value = 300000001
$ws.Range("G50").Value = "Rabanito"
$ws.Range("H50").Value = "Sin especificar"
$ws.Range("I50").Value = "Primera"
$ws.Range("J50").Value = 7900
$ws.Range("K50").Value = 3500
$ws.Range("L50").Value = 4000
$ws.Range("M50").Value = 3747
$ws.Range("N50").Value = '$/cien unidades (volumen en unidades)'
$ws.Range("O50").Value = "Provincia de Chacabuco"
$ws.Range("P50").Value = 37
$ws.Range("Q50").Value = 100
$ws.Range("R50").Value = "Hortaliza"
